# aula 15 do 8 ano
# Updates the "calendário" (schedule) sheet: the "Entrega" (delivery) dates
# in the 7º ano column (E) are rolled forward to the next cycle's dates,
# copying the corresponding values already present in the 6º ano column (D).
# The last activity (Atividade 4) had not yet been delivered, so its
# "Pedido"/"Entrega" cells (D40:D41 / E40:E41) are marked with a placeholder
# "-" just like column D already was.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("calendário")

# --- Content updates: column E (7 ano) rows 34-41 -------------------------
$ws.Range("E34").Value = 45565
$ws.Range("E35").Value = 45579
$ws.Range("E36").Value = 45579
$ws.Range("E37").Value = 45593
$ws.Range("E38").Value = 45593
$ws.Range("E39").Value = 45607

# Rows 40/41 (Atividade 4) mirror column D, which already holds a
# quote-prefixed "-" placeholder instead of a date.
$ws.Range("E40").Value = "'-"
$ws.Range("E41").Value = "'-"

# --- View bookkeeping: bring the new rows into view ------------------------
$ws.Activate()
$ws.Range("I47").Select()
$excel.ActiveWindow.Zoom = 190
